$d = $word.ActiveDocument

$pairs = @(
    @("490×5=", "385×3="),
    @("512×5=", "647×7="),
    @("155×6=", "605×8="),
    @("864×5=", "141×8="),
    @("395×4=", "837×2="),
    @("422×6=", "541×9="),
    @("520×6=", "532×9="),
    @("432×6=", "171×7="),
    @("895×9=", "169×6="),
    @("825×6=", "333×3="),
    @("646×5=", "182×4="),
    @("151×5=", "143×3="),
    @("888×2=", "620×6="),
    @("938×4=", "502×8="),
    @("839×7=", "317×9="),
    @("189×8=", "348×7="),
    @("160×7=", "519×6="),
    @("650×9=", "375×4="),
    @("683×7=", "469×7="),
    @("254×3=", "595×6="),
    @("114×5=", "609×6="),
    @("348×4=", "889×8="),
    @("587×4=", "969×9="),
    @("810×8=", "754×6="),
    @("752×6=", "384×8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
